$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1887640449438202
$ws.Range("C2").Value = 0.5595505617977528
$ws.Range("J2").Value = 0.01573033707865169
$ws.Range("P2").Value = 0.1460674157303371
$ws.Range("S2").Value = 0.0898876404494382

# Row 3
$ws.Range("B3").Value = 0.01587301587301587
$ws.Range("C3").Value = 0.007936507936507936
$ws.Range("J3").Value = 0.0119047619047619
$ws.Range("P3").Value = 0.7023809523809523
$ws.Range("S3").Value = 0.2619047619047619

# Row 5
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333

# Row 6
$ws.Range("B6").Value = 0.05639097744360902
$ws.Range("D6").Value = 0.01879699248120301
$ws.Range("F6").Value = 0.05263157894736842
$ws.Range("J6").Value = 0.2218045112781955
$ws.Range("O6").Value = 0.01503759398496241
$ws.Range("Q6").Value = 0.1240601503759398
$ws.Range("R6").Value = 0.06766917293233082
$ws.Range("S6").Value = 0.443609022556391

# Row 7
$ws.Range("B7").Value = 0.1290322580645161
$ws.Range("D7").Value = 0.02016129032258064
$ws.Range("F7").Value = 0.05241935483870968
$ws.Range("J7").Value = 0.0846774193548387
$ws.Range("O7").Value = 0.01612903225806452
$ws.Range("Q7").Value = 0.1451612903225807
$ws.Range("R7").Value = 0.06048387096774194
$ws.Range("S7").Value = 0.4919354838709677

# Row 8
$ws.Range("B8").Value = 0.1280788177339902
$ws.Range("D8").Value = 0.01642036124794746
$ws.Range("E8").Value = 0.004926108374384237
$ws.Range("F8").Value = 0.06075533661740559
$ws.Range("J8").Value = 0.09195402298850575
$ws.Range("O8").Value = 0.01642036124794746
$ws.Range("Q8").Value = 0.1412151067323481
$ws.Range("R8").Value = 0.07060755336617405
$ws.Range("S8").Value = 0.4696223316912972

# Row 9
$ws.Range("B9").Value = 0.155688622754491
$ws.Range("D9").Value = 0.02994011976047904
$ws.Range("F9").Value = 0.0718562874251497
$ws.Range("J9").Value = 0.1017964071856287
$ws.Range("O9").Value = 0.005988023952095809
$ws.Range("Q9").Value = 0.1497005988023952
$ws.Range("R9").Value = 0.04790419161676647
$ws.Range("S9").Value = 0.437125748502994

# Row 10
$ws.Range("B10").Value = 0.1479591836734694
$ws.Range("D10").Value = 0.02259475218658892
$ws.Range("F10").Value = 0.08309037900874636
$ws.Range("J10").Value = 0.1056851311953353
$ws.Range("O10").Value = 0.01457725947521866
$ws.Range("Q10").Value = 0.1712827988338192
$ws.Range("R10").Value = 0.06778425655976676
$ws.Range("S10").Value = 0.3870262390670554

# Row 11
$ws.Range("G11").Value = 0.145
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.2125
$ws.Range("L11").Value = 0.5275
$ws.Range("S11").Value = 0.015

# Row 12
$ws.Range("G12").Value = 0.7174887892376681
$ws.Range("J12").Value = 0.179372197309417
$ws.Range("K12").Value = 0.008968609865470852
$ws.Range("L12").Value = 0.02242152466367713
$ws.Range("S12").Value = 0.07174887892376682

# Row 13
$ws.Range("G13").Value = 0.7450980392156863
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.0196078431372549

# Row 15
$ws.Range("F15").Value = 0.03404255319148936
$ws.Range("H15").Value = 0.148936170212766
$ws.Range("I15").Value = 0.07659574468085106
$ws.Range("J15").Value = 0.3617021276595745
$ws.Range("K15").Value = 0.05957446808510639
$ws.Range("M15").Value = 0.008510638297872341
$ws.Range("N15").Value = 0.00425531914893617
$ws.Range("O15").Value = 0.05106382978723404
$ws.Range("S15").Value = 0.2553191489361702

# Row 16
$ws.Range("F16").Value = 0.01818181818181818
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.06545454545454546
$ws.Range("J16").Value = 0.3745454545454546
$ws.Range("K16").Value = 0.1163636363636364
$ws.Range("M16").Value = 0.01818181818181818
$ws.Range("O16").Value = 0.05818181818181818
$ws.Range("S16").Value = 0.1672727272727273

# Row 17
$ws.Range("F17").Value = 0.01445783132530121
$ws.Range("H17").Value = 0.2289156626506024
$ws.Range("I17").Value = 0.06024096385542169
$ws.Range("J17").Value = 0.3951807228915662
$ws.Range("K17").Value = 0.09156626506024096
$ws.Range("M17").Value = 0.02409638554216868
$ws.Range("O17").Value = 0.03132530120481928
$ws.Range("S17").Value = 0.1542168674698795

# Row 18
$ws.Range("F18").Value = 0.02840909090909091
$ws.Range("H18").Value = 0.2159090909090909
$ws.Range("I18").Value = 0.03977272727272727
$ws.Range("J18").Value = 0.4090909090909091
$ws.Range("K18").Value = 0.125
$ws.Range("O18").Value = 0.07954545454545454
$ws.Range("S18").Value = 0.1022727272727273

# Row 19
$ws.Range("F19").Value = 0.01195457262402869
$ws.Range("H19").Value = 0.2331141661685595
$ws.Range("I19").Value = 0.05977286312014345
$ws.Range("J19").Value = 0.3395098625224148
$ws.Range("K19").Value = 0.1243275552898984
$ws.Range("M19").Value = 0.02151823072325165
$ws.Range("N19").Value = 0.0005977286312014345
$ws.Range("O19").Value = 0.06814106395696354
$ws.Range("S19").Value = 0.1410639569635385
